# Added 4wk low sales check.
# - "Forecast Comparison": the Inventory Coverage (H) column is blanked
#   out for weeks where it no longer applies, and the Seasonality Index
#   (L) column is recalculated with updated values.
# - "Summary": the 4/8/16-week forecast totals are reset to 0 now that
#   the low-sales check zeroes out the short-term forecast.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Clear the "Inventory Coverage" (H) column values for rows 2-13 and 15
$ws1.Range("H2").Value = ""
$ws1.Range("H3").Value = ""
$ws1.Range("H4").Value = ""
$ws1.Range("H5").Value = ""
$ws1.Range("H6").Value = ""
$ws1.Range("H7").Value = ""
$ws1.Range("H8").Value = ""
$ws1.Range("H9").Value = ""
$ws1.Range("H10").Value = ""
$ws1.Range("H11").Value = ""
$ws1.Range("H12").Value = ""
$ws1.Range("H13").Value = ""
$ws1.Range("H15").Value = ""

# Update Seasonality Index (L) column values
$ws1.Range("L2").Value = 1
$ws1.Range("L3").Value = 0.96
$ws1.Range("L4").Value = 0.82
$ws1.Range("L5").Value = 0.92
$ws1.Range("L7").Value = 0.99
$ws1.Range("L8").Value = 1.01
$ws1.Range("L9").Value = 1.08
$ws1.Range("L11").Value = 1.09
$ws1.Range("L12").Value = 1.01
$ws1.Range("L13").Value = 0.82
$ws1.Range("L14").Value = 0.9399999999999999
$ws1.Range("L15").Value = 1.02
$ws1.Range("L16").Value = 0.88
$ws1.Range("L17").Value = 0.91

# Update Summary sheet totals to reflect the 4wk low sales check
# (force text formatting so the values are stored as strings, matching
#  the original inlineStr cell type rather than being auto-detected as numbers)
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "0"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "0"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "0"
